{"js": "// Docx writer: use a different style (\"Footnote Block Text\") for block\n// quotes that live inside footnotes, based on \"Footnote Text\" instead of\n// \"Block Text\"'s base (\"Body Text\"), so footnote block quotes can later be\n// given a different font size than footnote text.\n\n// 1. Mint the new paragraph style. (addStyle's own returned proxy can't be\n//    written to reliably before it is re-fetched through the styles\n//    collection, so sync once and then look it up by name.)\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nconst footnoteBlockText = styles.items.find(\n  (s) => s.nameLocal === \"Footnote Block Text\"\n);\nif (!footnoteBlockText) {\n  throw new Error('Could not find newly added style \"Footnote Block Text\"');\n}\n\n// 2. Base it on / chain it to \"Footnote Text\" (mirrors how \"Block Text\" is\n//    based on & followed by \"Body Text\").\nfootnoteBlockText.baseStyle = \"Footnote Text\";\nfootnoteBlockText.nextParagraphStyle = \"Footnote Text\";\n\n// 3. Same bookkeeping flags \"Block Text\" carries.\nfootnoteBlockText.priority = 9;\nfootnoteBlockText.unhideWhenUsed = true;\nfootnoteBlockText.quickStyle = true;\nawait context.sync();\n\n// 4. Paragraph formatting identical to \"Block Text\": 100 twips (5pt)\n//    spacing before/after, 480 twips (24pt) left/right indent, no\n//    first-line indent.\nconst pf = footnoteBlockText.paragraphFormat;\npf.spaceBefore = 5;\npf.spaceAfter = 5;\npf.leftIndent = 24;\npf.rightIndent = 24;\npf.firstLineIndent = 0;\nawait context.sync();\n", "ps1": "# Docx writer: use a different style (\"Footnote Block Text\") for block\n# quotes that live inside footnotes, based on \"Footnote Text\" instead of\n# \"Block Text\"'s base (\"Body Text\"), so footnote block quotes can later be\n# given a different font size than footnote text.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$style = $d.Styles.Add(\"Footnote Block Text\", 1)\n\n# Base it on / chain it to \"Footnote Text\" (mirrors how \"Block Text\" is\n# based on & followed by \"Body Text\").\n$style.BaseStyle = \"Footnote Text\"\n$style.NextParagraphStyle = \"Footnote Text\"\n\n# Same bookkeeping flags \"Block Text\" carries.\n$style.Priority = 9\n$style.UnhideWhenUsed = $true\n$style.QuickStyle = $true\n\n# Paragraph formatting identical to \"Block Text\": 100 twips (5pt) spacing\n# before/after, 480 twips (24pt) left/right indent, no first-line indent.\n$style.ParagraphFormat.SpaceBefore = 5\n$style.ParagraphFormat.SpaceAfter = 5\n$style.ParagraphFormat.LeftIndent = 24\n$style.ParagraphFormat.RightIndent = 24\n$style.ParagraphFormat.FirstLineIndent = 0\n"}
